$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add five new "Delegate" columns (O:S) -------------------------------
# 1) Write all the new values first.
$ws.Range("O1").Value = "Delegate First Name"
$ws.Range("P1").Value = "Delegate Last Name"
$ws.Range("Q1").Value = "Delegate Email"
$ws.Range("R1").Value = "Delegate Phone Number"
$ws.Range("S1").Value = "Delegate Phone Type"

$ws.Range("O2").Value = "Text capitalized as you want the name to appear in the platform."
$ws.Range("P2").Value = "Text"
$ws.Range("Q2").Value = "youremail@yourdomain.com"
$ws.Range("R2").Value = "Format per country 919-555-1212"
$ws.Range("S2").Value = "Enter one of these values: Home        Work      Mobile"

# 2) Now stamp each new cell with the formatting of its sibling column so the
#    workbook keeps reusing the existing cell-style table instead of growing
#    new (duplicate) styles.
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1:S1").PasteSpecial(-4122) | Out-Null

$ws.Range("N2").Copy() | Out-Null
$ws.Range("O2").PasteSpecial(-4122) | Out-Null
$ws.Range("R2").PasteSpecial(-4122) | Out-Null
$ws.Range("S2").PasteSpecial(-4122) | Out-Null

$ws.Range("M2").Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null

$ws.Range("G2").Copy() | Out-Null
$ws.Range("Q2").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Select() | Out-Null

# 3) Mail-to hyperlink on the new "Delegate Email" example cell, matching
#    the existing one on G2.
$ws.Hyperlinks.Add($ws.Range("Q2"), "mailto:youremail@yourdomain.com") | Out-Null

# --- View state: zoom + selection ----------------------------------------
$excel.ActiveWindow.Zoom = 60
$ws.Range("S1").Select() | Out-Null

Write-Output "done"
